$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.168.65'
$ws.Range("E2").Value = '  -1.00%  '

$ws.Range("D3").Value = '1.852.22'
$ws.Range("E3").Value = '  -2.35%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.63'
$ws.Range("E5").Value = '  -1.47%  '

$ws.Range("E6").Value = '  +0.16%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4758'
$ws.Range("E7").Value = '  -2.95%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2785'
$ws.Range("E8").Value = '  -4.84%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06387'
$ws.Range("E9").Value = '  -4.58%  '

$ws.Range("D10").Value = '1.855.12'
$ws.Range("E10").Value = '  -2.30%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07358'
$ws.Range("E11").Value = '  +0.36%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.03'
$ws.Range("E12").Value = '  -5.48%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.070'
$ws.Range("E13").Value = '  -2.00%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '86.39'
$ws.Range("E14").Value = '  -1.51%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6404'
$ws.Range("E15").Value = '  -3.85%  '

$ws.Range("D16").Value = '30.112.68'
$ws.Range("E16").Value = '  -1.07%  '

$ws.Range("E17").Value = '  +0.03%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.05'
$ws.Range("E18").Value = '  -3.28%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007495'
$ws.Range("E19").Value = '  -4.87%  '

$ws.Range("B20").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C20").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D20").Value = '2.098.43'
$ws.Range("E20").Value = '  -2.35%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  +0.08%  '

$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.245'
$ws.Range("E22").Value = '  -1.61%  '

$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '220.02'
$ws.Range("E23").Value = '  +13.52%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.049'
$ws.Range("E24").Value = '  -1.24%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.154'
$ws.Range("E25").Value = '  -3.63%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.87'
$ws.Range("E26").Value = '  +0.09%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.22'
$ws.Range("E27").Value = '  -0.34%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.901'
$ws.Range("E28").Value = '  -1.98%  '

$ws.Range("E29").Value = '  -2.22%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09127'
$ws.Range("E30").Value = '  -0.40%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.178'
$ws.Range("E31").Value = '  -3.59%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.919'
$ws.Range("E32").Value = '  -3.60%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04938'
$ws.Range("E33").Value = '  -4.32%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.132'
$ws.Range("E34").Value = '  +2.69%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7160'
$ws.Range("E35").Value = '  -3.39%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.692'

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01819'
$ws.Range("E37").Value = '  +0.54%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.594'
$ws.Range("E38").Value = '  -3.54%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8970'
$ws.Range("E39").Value = '  -2.89%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.015'
$ws.Range("E40").Value = '  -1.56%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.855'
$ws.Range("E41").Value = '  -1.08%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '105.06'
$ws.Range("E42").Value = '  -1.70%  '

$ws.Range("E43").Value = '  +0.67%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4210'
$ws.Range("E44").Value = '  -4.00%  '

$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1297'
$ws.Range("E45").Value = '  -5.43%  '

$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.224'
$ws.Range("E46").Value = '  -4.66%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '63.47'
$ws.Range("E47").Value = '  -8.08%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.486'
$ws.Range("E48").Value = '  +5.84%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.663'
$ws.Range("E49").Value = '  -3.80%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '33.40'
$ws.Range("E50").Value = '  -4.38%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05642'
$ws.Range("E51").Value = '  -3.54%  '
